# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the "vintage" date for each forecast row. It was being
# written as the 1st of the vintage's reference month, but the correct
# convention is the 15th of the *following* month. Shift every date in
# column A from <month, day 1> -> <next month, day 15>, leaving every
# other value/style in the sheet untouched.

function Get-ExcelSerial($y, $m, $d) {
    # Gregorian calendar date -> Julian Day Number -> Excel 1900-date-system
    # serial (serial 0 == 1899-12-30), via the standard civil-to-JDN formula.
    if ($m -le 2) {
        $y = $y - 1
        $m = $m + 12
    }
    $a = [Math]::Floor((14 - $m) / 12)
    $y2 = $y + 4800 - $a
    $m2 = $m + 12 * $a - 3
    $jdn = $d + [Math]::Floor((153 * $m2 + 2) / 5) + 365 * $y2 + [Math]::Floor($y2 / 4) - [Math]::Floor($y2 / 100) + [Math]::Floor($y2 / 400) - 32045
    return $jdn - 2415019
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null -or $serial -eq "") {
        continue
    }

    $year = $excel.WorksheetFunction.Year($serial)
    $month = $excel.WorksheetFunction.Month($serial)

    $newMonth = $month + 1
    $newYear = $year
    if ($newMonth -gt 12) {
        $newMonth = $newMonth - 12
        $newYear = $newYear + 1
    }

    $newSerial = Get-ExcelSerial $newYear $newMonth 15
    $cell.Value = $newSerial
}
